$wb = $excel.ActiveWorkbook

# --- "Customers" sheet -----------------------------------------------------
$wsCustomers = $wb.Worksheets.Item("Customers")

# The "Test Result" boolean in F2 was wrong - flip it from TRUE to FALSE.
$wsCustomers.Range("F2").Value = $false

# The rest of the "Test Result" column (F6:F14) was bogus test data that
# should not have been there at all - clear it out completely.
$wsCustomers.Range("F6:F14").ClearContents()

# --- "Orders" sheet ---------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Orders")
[void]$wsOrders.Range("E2").Select()

# Restore "Customers" as the active sheet / tab, with its original selection.
$wsCustomers.Activate()
[void]$wsCustomers.Range("F2").Select()
